# Automatische test-sync: 2025-07-23 22:27:50
#
# Appends a new "Testmail #6" row to the Logs sheet, updates the matching
# category tally on the Dashboard sheet, and extends the bar chart's
# category/value series references to include the new Dashboard row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Logs sheet: append row 16
# ---------------------------------------------------------------
$wsLogs = $wb.Worksheets.Item("Logs")

$wsLogs.Range("A16").Value = "Zou jij deze offerte even naar klant X willen mailen?"
$wsLogs.Range("B16").Value = "mailmind.test@zohomail.eu"
$wsLogs.Range("C16").Value = "Testmail #6: Zou jij deze offerte even naar klant X willen mailen?"
$wsLogs.Range("D16").Value = "Offerte / Prijsaanvraag"
$wsLogs.Range("E16").Value = "Geachte afzender,`nBedankt voor uw e-mail. Helaas kan ik u niet helpen bij het versturen van de offerte naar klant X, aangezien ik geen toegang heb tot uw documenten of het vermogen om e-mails namens u te verzenden.`nIk raad u aan de offerte zelf naar klant X te mailen of contact op te nemen met de verantwoordelijke persoon binnen uw organisatie die dit kan afhandelen.`nAls u nog andere vragen heeft, laat het me dan weten.`nMet vriendelijke groet,`n[Naam]`nE-mailassistent"
$wsLogs.Range("F16").Value = "2025-07-23 22:26:57"
$wsLogs.Range("G16").Value = "Ja"
$wsLogs.Range("H16").Value = "Nee"
$wsLogs.Range("I16").Value = "Ja"
$wsLogs.Range("J16").Value = "Nee"

# ---------------------------------------------------------------
# 2. Logs sheet: extend the conditional-formatting ranges from row 15
#    to row 16 so the newly added row is covered too.
# ---------------------------------------------------------------
$wsLogs.Range("D2:D15").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("D2:D16"))
$wsLogs.Range("G2:G15").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("G2:G16"))
$wsLogs.Range("H2:H15").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("H2:H16"))
$wsLogs.Range("I2:I15").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("I2:I16"))
$wsLogs.Range("J2:J15").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("J2:J16"))

# ---------------------------------------------------------------
# 3. Dashboard sheet: append tally row 7 for the new category
# ---------------------------------------------------------------
$wsDash = $wb.Worksheets.Item("Dashboard")

$wsDash.Range("A7").Value = "Offerte / Prijsaanvraag"
$wsDash.Range("B7").Value = 1

# ---------------------------------------------------------------
# 4. Dashboard sheet: extend the bar chart's category/value series
#    so it now spans rows 2-7 instead of 2-6.
# ---------------------------------------------------------------
$chart = $wsDash.ChartObjects(1).Chart
$ser = $chart.SeriesCollection(1)
$ser.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$7,Dashboard!`$B`$2:`$B`$7,1)"
